$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2 through 11
# from serial 46070 (2026-02-17) to serial 46072 (2026-02-19)
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46070) {
        $cell.Value2 = 46072
    }
}
